$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Update rows 2-8 and 10-15 in place (a permutation of existing records' data).
# Row 9 is untouched. Only columns A,B,D,E,F,G,H,Q,R change values; row 13/14 also
# exchange the 'Publik kommentar' (AC) note that travels together with the moved record.

# Row 2: set to record previously at a different row (Id=111523728)
$ws.Range("A2").Value = 111523728
$ws.Range("B2").Value = 89686
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 658
$ws.Range("F2").Value = 'Rosenticka'
$ws.Range("G2").Value = 'Rhodofomes roseus'
$ws.Range("H2").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q2").Value = 497338.5868253836
$ws.Range("R2").Value = 6754122.194367126

# Row 3: set to record previously at a different row (Id=111523701)
$ws.Range("A3").Value = 111523701
$ws.Range("B3").Value = 89686
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 658
$ws.Range("F3").Value = 'Rosenticka'
$ws.Range("G3").Value = 'Rhodofomes roseus'
$ws.Range("H3").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q3").Value = 497367.2942720717
$ws.Range("R3").Value = 6754083.757028132

# Row 4: set to record previously at a different row (Id=111523727)
$ws.Range("A4").Value = 111523727
$ws.Range("B4").Value = 89845
$ws.Range("D4").Value = 'VU'
$ws.Range("E4").Value = 1209
$ws.Range("F4").Value = 'Rynkskinn'
$ws.Range("G4").Value = 'Phlebia centrifuga'
$ws.Range("H4").Value = 'P.Karst.'
$ws.Range("Q4").Value = 497338.5868253836
$ws.Range("R4").Value = 6754122.194367126

# Row 5: set to record previously at a different row (Id=111523712)
$ws.Range("A5").Value = 111523712
$ws.Range("B5").Value = 89405
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 1202
$ws.Range("F5").Value = 'Ullticka'
$ws.Range("G5").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H5").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q5").Value = 497301.0581945881
$ws.Range("R5").Value = 6754088.183226726

# Row 6: set to record previously at a different row (Id=111523730)
$ws.Range("A6").Value = 111523730
$ws.Range("B6").Value = 89405
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 1202
$ws.Range("F6").Value = 'Ullticka'
$ws.Range("G6").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H6").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q6").Value = 497338.5868253836
$ws.Range("R6").Value = 6754122.194367126

# Row 7: set to record previously at a different row (Id=111523741)
$ws.Range("A7").Value = 111523741
$ws.Range("B7").Value = 89686
$ws.Range("D7").Value = 'NT'
$ws.Range("E7").Value = 658
$ws.Range("F7").Value = 'Rosenticka'
$ws.Range("G7").Value = 'Rhodofomes roseus'
$ws.Range("H7").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q7").Value = 497384.3941364431
$ws.Range("R7").Value = 6754155.713205664

# Row 8: set to record previously at a different row (Id=111523683)
$ws.Range("A8").Value = 111523683
$ws.Range("B8").Value = 89845
$ws.Range("D8").Value = 'VU'
$ws.Range("E8").Value = 1209
$ws.Range("F8").Value = 'Rynkskinn'
$ws.Range("G8").Value = 'Phlebia centrifuga'
$ws.Range("H8").Value = 'P.Karst.'
$ws.Range("Q8").Value = 497391.6869587752
$ws.Range("R8").Value = 6754138.20205555

# Row 10: set to record previously at a different row (Id=111523657)
$ws.Range("A10").Value = 111523657
$ws.Range("B10").Value = 89686
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 658
$ws.Range("F10").Value = 'Rosenticka'
$ws.Range("G10").Value = 'Rhodofomes roseus'
$ws.Range("H10").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q10").Value = 497390.1961838813
$ws.Range("R10").Value = 6754097.842248607

# Row 11: set to record previously at a different row (Id=111523731)
$ws.Range("A11").Value = 111523731
$ws.Range("B11").Value = 89845
$ws.Range("D11").Value = 'VU'
$ws.Range("E11").Value = 1209
$ws.Range("F11").Value = 'Rynkskinn'
$ws.Range("G11").Value = 'Phlebia centrifuga'
$ws.Range("H11").Value = 'P.Karst.'
$ws.Range("Q11").Value = 497307.3714758331
$ws.Range("R11").Value = 6754063.864355386

# Row 12: set to record previously at a different row (Id=111523656)
$ws.Range("A12").Value = 111523656
$ws.Range("B12").Value = 89845
$ws.Range("D12").Value = 'VU'
$ws.Range("E12").Value = 1209
$ws.Range("F12").Value = 'Rynkskinn'
$ws.Range("G12").Value = 'Phlebia centrifuga'
$ws.Range("H12").Value = 'P.Karst.'
$ws.Range("Q12").Value = 497390.1961838813
$ws.Range("R12").Value = 6754097.842248607

# Row 13: set to record previously at a different row (Id=111523697)
$ws.Range("A13").Value = 111523697
$ws.Range("B13").Value = 77515
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = 'Garnlav'
$ws.Range("G13").Value = 'Alectoria sarmentosa'
$ws.Range("H13").Value = '(Ach.) Ach.'
$ws.Range("Q13").Value = 497380.5053056676
$ws.Range("R13").Value = 6754165.927741241

# Row 14: set to record previously at a different row (Id=111523740)
$ws.Range("A14").Value = 111523740
$ws.Range("B14").Value = 56398
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 100109
$ws.Range("F14").Value = 'Tretåig hackspett'
$ws.Range("G14").Value = 'Picoides tridactylus'
$ws.Range("H14").Value = '(Linnaeus, 1758)'
$ws.Range("Q14").Value = 497366.3615979423
$ws.Range("R14").Value = 6754139.679549156

# Row 15: set to record previously at a different row (Id=111523695)
$ws.Range("A15").Value = 111523695
$ws.Range("B15").Value = 5113
$ws.Range("D15").Value = 'LC'
$ws.Range("E15").Value = 100526
$ws.Range("F15").Value = 'Bronshjon'
$ws.Range("G15").Value = 'Callidium coriaceum'
$ws.Range("H15").Value = 'Paykull, 1800'
$ws.Range("Q15").Value = 497354.1644349985
$ws.Range("R15").Value = 6754111.484663551

# Move the 'Publik kommentar' note from row 13 to row 14 (it belongs to Id 111523740,
# which is now located in row 14).
$ws.Range("AC14").Value = 'Gamla hack'
$ws.Range("AC13").Value = $null
